$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1100, 3.951002359390259, 1.281787157058716, -0.6831467747688293, 9.308746337890625, 18476, "01"),
    @(1200, 3.904142141342163, 1.283282518386841, -0.9052648544311523, 10.12545680999756, 18568, "01"),
    @(1300, 3.842416048049927, 1.282211780548096, -0.9531055688858032, 10.13229084014893, 18438, "01"),
    @(1400, 3.818761587142944, 1.266582608222961, -0.8215433955192566, 8.560379028320312, 18366, "01"),
    @(1500, 3.787439107894897, 1.304946541786194, -1.019741058349609, 10.73371887207031, 18392, "01")
)

$startRow = 21
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $rowData[6]
    $gCell.NumberFormat = "General"
    $gCell.Style = "Normal"
}
